# Adds a new "2022-Q1" fund-holdings sheet (cloned structure/styles from
# "2021-Q4") right before "总计", and inserts the matching 2022-Q1 summary
# row at the top of "总计".

$wb = $excel.ActiveWorkbook

# ---------- 1) new "2022-Q1" detail sheet ----------
$srcDetail = $wb.Worksheets.Item("2021-Q4")

# Insert directly after "2021-Q4" (i.e. right before "总计").
$newSheet = $wb.Worksheets.Add($null, $srcDetail)
$newSheet.Name = "2022-Q1"

# Clone the full A1:H18 layout (headers, borders, bold/centered styles,
# column A running index) from the previous quarter sheet...
$srcDetail.Range("A1:H18").Copy($newSheet.Range("A1:H18"))
$newSheet.Range("A1").ClearContents()

# ...then overwrite with 2022-Q1 figures.
#
# Columns B:G hold figures like fund code "005680" or size "57.61" that must
# stay literal TEXT (as in the source data), not auto-coerced numbers/leading
# zeros dropped -- briefly mark the range as Text before writing, then reset
# the style back to Normal (these cells carry no special formatting of their
# own, only the string values) so only the A/H numeric columns and the header
# row keep the bold/bordered style copied above.
$dataRows = $newSheet.Range("B2:G18")
$dataRows.NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = '005680'
$newSheet.Range("C2").Value = '财通资管价值成长混合'
$newSheet.Range("D2").Value = '57.61'
$newSheet.Range("E2").Value = '91.74'
$newSheet.Range("F2").Value = '4.41'
$newSheet.Range("G2").Value = '2.5406'
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = '010846'
$newSheet.Range("C3").Value = '南方卓越优选3个月持有期混合A'
$newSheet.Range("D3").Value = '26.01'
$newSheet.Range("E3").Value = '60.36'
$newSheet.Range("F3").Value = '5.56'
$newSheet.Range("G3").Value = '1.4462'
$newSheet.Range("H3").Value = 3

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = '008276'
$newSheet.Range("C4").Value = '财通资管价值发现混合'
$newSheet.Range("D4").Value = '21.19'
$newSheet.Range("E4").Value = '92.22'
$newSheet.Range("F4").Value = '4.78'
$newSheet.Range("G4").Value = '1.0129'
$newSheet.Range("H4").Value = 8

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = '009950'
$newSheet.Range("C5").Value = '财通资管均衡价值一年持有期混合'
$newSheet.Range("D5").Value = '21.22'
$newSheet.Range("E5").Value = '91.17'
$newSheet.Range("F5").Value = '4.08'
$newSheet.Range("G5").Value = '0.8658'
$newSheet.Range("H5").Value = 9

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = '010163'
$newSheet.Range("C6").Value = '财通资管价值精选一年持有期混合A'
$newSheet.Range("D6").Value = '16.79'
$newSheet.Range("E6").Value = '90.42'
$newSheet.Range("F6").Value = '4.05'
$newSheet.Range("G6").Value = '0.6800'
$newSheet.Range("H6").Value = 9

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = '010413'
$newSheet.Range("C7").Value = '财通资管宸瑞一年持有期混合A'
$newSheet.Range("D7").Value = '12.66'
$newSheet.Range("E7").Value = '94.71'
$newSheet.Range("F7").Value = '5.17'
$newSheet.Range("G7").Value = '0.6545'
$newSheet.Range("H7").Value = 8

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = '011488'
$newSheet.Range("C8").Value = '申万菱信乐享混合'
$newSheet.Range("D8").Value = '12.18'
$newSheet.Range("E8").Value = '77.88'
$newSheet.Range("F8").Value = '3.70'
$newSheet.Range("G8").Value = '0.4507'
$newSheet.Range("H8").Value = 7

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = '013634'
$newSheet.Range("C9").Value = '申万菱信双利混合A'
$newSheet.Range("D9").Value = '7.83'
$newSheet.Range("E9").Value = '22.26'
$newSheet.Range("F9").Value = '2.93'
$newSheet.Range("G9").Value = '0.2294'
$newSheet.Range("H9").Value = 3

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = '010847'
$newSheet.Range("C10").Value = '南方卓越优选3个月持有期混合C'
$newSheet.Range("D10").Value = '3.42'
$newSheet.Range("E10").Value = '60.36'
$newSheet.Range("F10").Value = '5.56'
$newSheet.Range("G10").Value = '0.1902'
$newSheet.Range("H10").Value = 3

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = '010414'
$newSheet.Range("C11").Value = '财通资管宸瑞一年持有期混合C'
$newSheet.Range("D11").Value = '1.53'
$newSheet.Range("E11").Value = '94.71'
$newSheet.Range("F11").Value = '5.17'
$newSheet.Range("G11").Value = '0.0791'
$newSheet.Range("H11").Value = 8

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = '004703'
$newSheet.Range("C12").Value = '南方兴盛先锋灵活配置混合'
$newSheet.Range("D12").Value = '1.09'
$newSheet.Range("E12").Value = '53.08'
$newSheet.Range("F12").Value = '3.82'
$newSheet.Range("G12").Value = '0.0416'
$newSheet.Range("H12").Value = 4

$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = '310368'
$newSheet.Range("C13").Value = '申万菱信竞争优势混合'
$newSheet.Range("D13").Value = '0.83'
$newSheet.Range("E13").Value = '91.22'
$newSheet.Range("F13").Value = '4.40'
$newSheet.Range("G13").Value = '0.0365'
$newSheet.Range("H13").Value = 5

$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = '005444'
$newSheet.Range("C14").Value = '光大保德信多策略精选18个月定期开放灵活配置混合'
$newSheet.Range("D14").Value = '1.00'
$newSheet.Range("E14").Value = '29.09'
$newSheet.Range("F14").Value = '2.59'
$newSheet.Range("G14").Value = '0.0259'
$newSheet.Range("H14").Value = 4

$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = '013635'
$newSheet.Range("C15").Value = '申万菱信双利混合C'
$newSheet.Range("D15").Value = '0.75'
$newSheet.Range("E15").Value = '22.26'
$newSheet.Range("F15").Value = '2.93'
$newSheet.Range("G15").Value = '0.0220'
$newSheet.Range("H15").Value = 3

$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = '010164'
$newSheet.Range("C16").Value = '财通资管价值精选一年持有期混合C'
$newSheet.Range("D16").Value = '0.47'
$newSheet.Range("E16").Value = '90.42'
$newSheet.Range("F16").Value = '4.05'
$newSheet.Range("G16").Value = '0.0190'
$newSheet.Range("H16").Value = 9

$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = '001148'
$newSheet.Range("C17").Value = '申万菱信多策略灵活配置混合A'
$newSheet.Range("D17").Value = '1.10'
$newSheet.Range("E17").Value = '28.34'
$newSheet.Range("F17").Value = '1.45'
$newSheet.Range("G17").Value = '0.0160'
$newSheet.Range("H17").Value = 4

$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = '001724'
$newSheet.Range("C18").Value = '申万菱信多策略灵活配置混合C'
$newSheet.Range("D18").Value = '1.06'
$newSheet.Range("E18").Value = '28.34'
$newSheet.Range("F18").Value = '1.45'
$newSheet.Range("G18").Value = '0.0154'
$newSheet.Range("H18").Value = 4

# Drop the temporary Text number format back to Normal (style index 0,
# matching the source data rows) now that the literal strings are in place.
$dataRows.Style = "Normal"

# ---------- 2) prepend 2022-Q1 row into "总计" ----------
$totalSheet = $wb.Worksheets.Item('总计')

# Make room for the new row, formatted like a normal data row (copy format
# down from the row that is about to become row 3, then write the new data).
$totalSheet.Rows("2:2").Insert()
$totalSheet.Range("A2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = '2022-Q1'
$totalSheet.Range("C2").Value = 17
$totalSheet.Range("D2").Value = 8.33

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = '2021-Q4'
$totalSheet.Range("C3").Value = 17
$totalSheet.Range("D3").Value = 9.35

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = '2021-Q3'
$totalSheet.Range("C4").Value = 57
$totalSheet.Range("D4").Value = 30.51

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = '2021-Q2'
$totalSheet.Range("C5").Value = 46
$totalSheet.Range("D5").Value = 22.29

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = '2021-Q1'
$totalSheet.Range("C6").Value = 31
$totalSheet.Range("D6").Value = 13.56

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = '2020-Q4'
$totalSheet.Range("C7").Value = 21
$totalSheet.Range("D7").Value = 7.01
